$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The export template originally had 8 data columns (A..H): date, account,
# address, unitName, type, transaction, status, amount. This change adds a
# new "order" column between "transaction" (F/G) and "status" (old G/H),
# i.e. a brand-new column G, pushing status/amount from G/H to H/I.
#
# Insert a new column at G - this shifts old G (status) and H (amount)
# one column to the right (to H and I respectively), while copying the
# row-by-row formatting of the old column G onto both the new G and the
# shifted H, which matches how the source file's styles line up.
$ws.Range("G1").EntireColumn.Insert()

# Header row (row 1): i18n placeholder for the new "order" column.
$ws.Cells.Item(1, 7).Value = "{d.i18n.order}"

# Sample/template rows (rows 2 and 3): the data-binding placeholders for
# the new "order" field, first iterator and "next" iterator variants.
$ws.Cells.Item(2, 7).Value = "{d.objs[I].order}"
$ws.Cells.Item(3, 7).Value = "{d.objs[I+1].order}"

# Give the new "order" column its own (wider) width; the other columns
# keep their previous widths (H/I now share the old G/H width of 14.5).
$ws.Range("G1").EntireColumn.ColumnWidth = 25.9375
